{"js": "// Office.js (Word JavaScript API) edit script.\n// This is the body of `async (context) => { ... }`.\n//\n// The underlying edit (per the commit \"ajout de ma partie ?\") rewrites the\n// sentence describing the \"gestion des cong\u00e9s\" (leave/holiday management)\n// feature: the old, generic wording about adding date-range periods is\n// replaced with a more specific description about adding/removing a leave\n// day by day.\n//\n// Old: \"La gestion des cong\u00e9s permettra d\u2019ajouter des p\u00e9riodes (date d\u00e9but\n//       + date fin) de les modifier et de les supprimer.\"\n// New: \"La gestion des cong\u00e9s permettra de supprimer un cong\u00e9 ou\n//       d\u2019ajouter un cong\u00e9, jour par jour\"\n\nconst body = context.document.body;\n\nconst oldTail =\n  \"permettra d\\u2019ajouter des p\\u00e9riodes (date d\\u00e9but + date fin) de les modifier et de les supprimer.\";\nconst newTail =\n  \"permettra de supprimer un cong\\u00e9 ou d\\u2019ajouter un cong\\u00e9, jour par jour\";\n\nconst searchResults = body.search(oldTail, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  // Replace only the trailing part of the sentence, leaving the\n  // \"La gestion des cong\u00e9s \" lead-in (and its run/formatting) untouched.\n  searchResults.items[0].insertText(newTail, Word.InsertLocation.replace);\n} else {\n  // Fallback: if the exact phrase can't be found (e.g. text already\n  // changed), locate the paragraph by its distinctive lead-in text and\n  // replace its content wholesale, preserving paragraph formatting.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  const target = paragraphs.items.find((p) =>\n    p.text.indexOf(\"La gestion des cong\\u00e9s\") !== -1\n  );\n  if (target) {\n    target.insertText(\n      \"La gestion des cong\\u00e9s permettra de supprimer un cong\\u00e9 ou d\\u2019ajouter un cong\\u00e9, jour par jour\",\n      Word.InsertLocation.replace\n    );\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# The underlying edit (per the commit \"ajout de ma partie ?\") rewrites the\n# sentence describing the \"gestion des conges\" (leave/holiday management)\n# feature: the old, generic wording about adding date-range periods is\n# replaced with a more specific description about adding/removing a leave\n# day by day.\n#\n# Old: \"La gestion des conges permettra d'ajouter des periodes (date debut\n#       + date fin) de les modifier et de les supprimer.\"\n# New: \"La gestion des conges permettra de supprimer un conge ou\n#       d'ajouter un conge, jour par jour\"\n\n$d = $word.ActiveDocument\n\n$oldTail = \"permettra d\u2019ajouter des p\u00e9riodes (date d\u00e9but + date fin) de les modifier et de les supprimer.\"\n$newTail = \"permettra de supprimer un cong\u00e9 ou d\u2019ajouter un cong\u00e9, jour par jour\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldTail\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newTail\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback: locate the paragraph by its distinctive lead-in text and\n    # replace its content wholesale, preserving paragraph formatting.\n    $leadIn = \"La gestion des cong\u00e9s\"\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $para = $d.Paragraphs.Item($i)\n        if ($para.Range.Text.Contains($leadIn)) {\n            $r = $para.Range\n            $r.End = $r.End - 1\n            $r.Text = \"La gestion des cong\u00e9s permettra de supprimer un cong\u00e9 ou d\u2019ajouter un cong\u00e9, jour par jour\"\n            break\n        }\n    }\n}\n"}
